$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: A2 and E2 become numeric values instead of text
$ws.Range("A2").Value = 612203144
$ws.Range("E2").Value = 9673095937

# Row 3: new registration (numeric MIS / phone)
$ws.Range("A3").Value = 612203154
$ws.Range("B3").Value = "Sakshi"
$ws.Range("C3").Value = "Khanorkar"
$ws.Range("D3").Value = "sakshikhanorkar15@gmail.com"
$ws.Range("E3").Value = 1234567890

# Row 4: new registration (MIS / phone stored as text)
$ws.Range("A4").NumberFormat = "@"
$ws.Range("A4").Value = "612203142"
$ws.Range("A4").Style = "Normal"

$ws.Range("B4").Value = "Yash"
$ws.Range("C4").Value = "Pawar"
$ws.Range("D4").Value = "yashpawar123@gmail.com"

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "1234567890"
$ws.Range("E4").Style = "Normal"

Write-Output "applied edits"
